$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.386.80"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.848.52"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'240.34"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "'0.6277"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.07629"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").Value = "'0.2908"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").Value = "'24.68"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").Value = "'0.07731"
$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").Value = "'5.027"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").Value = "'0.6787"
$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("E14").Value = "  -2.44%  "

$ws.Range("D15").Value = "'83.18"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "'6.155"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "29.431.39"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'226.56"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "'12.32"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").Value = "'7.482"
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'157.91"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("E24").Value = "  -0.90%  "

$ws.Range("D25").Value = "'8.393"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D27").Value = "'1.387"
$ws.Range("E27").Value = "  +5.63%  "

$ws.Range("D28").Value = "'1.461"
$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("D29").Value = "'0.05602"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("E30").Value = "  +0.50%  "

$ws.Range("D31").Value = "'4.065"
$ws.Range("E31").Value = "  +0.69%  "

$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").Value = "'0.6950"
$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").Value = "'2.585"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "1.228.49"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").Value = "'2.722"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("D39").Value = "'6.384"
$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").Value = "'0.9030"
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").Value = "'101.65"

$ws.Range("D43").Value = "'65.92"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").Value = "'7.158"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("D45").Value = "'0.00000000116"
$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("D47").Value = "'9.014"
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").Value = "'1.681"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("E49").Value = "  +2.05%  "

$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").Value = "'0.4628"
$ws.Range("E51").Value = "  +0.08%  "
